$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("maze")

# "Bonus" square becomes "B", and the second "Bonus" cell becomes "BBO"
$ws.Range("C3").Value = "B"
$ws.Range("C4").Value = "BBO"

# Beans squares VBx renamed to VNx
$ws.Range("B8").Value = "VNU"
$ws.Range("C8").Value = "VNR"
$ws.Range("D8").Value = "VND"
$ws.Range("E8").Value = "VNL"

# Update the active selection on the "maze" sheet to E9
$ws.Range("E9").Select()
